$d = $word.ActiveDocument

# The edit removes every paragraph that follows the "Author" paragraph
# ("Sam Abbott") -- i.e. the "Learning Objectives" section through to the
# final "References" heading are all deleted, leaving just the Title and
# Author paragraphs (and the trailing section mark).

$authorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Author") {
        $authorIndex = $i
    }
}

if ($authorIndex -gt 0 -and $authorIndex -lt $d.Paragraphs.Count) {
    $authorPara = $d.Paragraphs.Item($authorIndex)
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

    $startOfDeletion = $authorPara.Range.End
    $endOfDeletion = $lastPara.Range.End

    $deleteRange = $d.Range($startOfDeletion, $endOfDeletion)
    $deleteRange.Delete()
}

Write-Output ("Paragraphs remaining: " + $d.Paragraphs.Count)
